$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Constant"/"r2_adj" row entirely (content + formatting).
$ws.Range("A5:C5").Clear()

# New "C/A" column header, copying the existing bold/boxed header style
# from B1 (re-using the workbook's existing style slot rather than
# re-deriving it one property at a time, which would otherwise leave
# orphaned cellXfs entries behind).
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Header row text ---
$ws.Range("C1").Value = "C/A"
$ws.Range("D1").Value = "FFR"

# --- Row labels (style already in place on these cells) ---
$ws.Range("A3").Value = "C/A Lag"
$ws.Range("A4").Value = "FFR Lag"

# --- Data values ---
$ws.Range("B2").Value = "-0.656***"
$ws.Range("C2").Value = "2.634***"

# "0.039" looks like a plain number to Excel's type-inference, so force
# text storage (matching the source file, which keeps it as a shared
# string): enter it as a text-returning formula, then paste-special the
# computed result back over itself as a value. This keeps the literal
# string without ever touching NumberFormat (which would otherwise mint
# a permanent, unused style slot in styles.xml).
$ws.Range("D2").Formula = '="0.039"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("B3").Value = "-0.016***"
$ws.Range("C3").Value = "-0.763***"
$ws.Range("D3").Value = "-0.033***"

$ws.Range("B4").Value = "0.352***"
$ws.Range("C4").Value = "9.463***"
$ws.Range("D4").Value = "0.424***"
